$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Apply the underline-font style first, so it becomes cellXfs index 3
# (matches the new <font><u/>...</font> entry + new cellXfs entry in the
# target styles.xml, which is created before the yellow-fill entry).
$ws.Range("P33").Font.Underline = $true

# --- New experiment rows (retrained 5verbs runs), rows 40-43 ---
$row40 = New-Object 'object[,]' 1,15
$row40[0,0]  = "IA-PUCP"
$row40[0,1]  = "5verbs - retrained"
$row40[0,2]  = 54
$row40[0,3]  = 10
$row40[0,4]  = "e-4"
$row40[0,5]  = 100
$row40[0,6]  = 1
$row40[0,7]  = 8
$row40[0,8]  = 8
$row40[0,9]  = 0
$row40[0,10] = 0.1188
$row40[0,11] = 1
$row40[0,12] = 0.41389999999999999
$row40[0,13] = 0.86960000000000004
$row40[0,14] = 0.21641791044776101
$ws.Range("A40:O40").Value = $row40

$row41 = New-Object 'object[,]' 1,15
$row41[0,0]  = "IA-PUCP"
$row41[0,1]  = "5verbs - retrained"
$row41[0,2]  = 55
$row41[0,3]  = 10
$row41[0,4]  = "e-5"
$row41[0,5]  = 100
$row41[0,6]  = 1
$row41[0,7]  = 8
$row41[0,8]  = 8
$row41[0,9]  = 0
$row41[0,10] = 1.1040000000000001
$row41[0,11] = 0.72629999999999995
$row41[0,12] = 1.2490000000000001
$row41[0,13] = 0.52170000000000005
$row41[0,14] = 0.57462686567164101
$ws.Range("A41:O41").Value = $row41

$row42 = New-Object 'object[,]' 1,15
$row42[0,0]  = "IA-PUCP"
$row42[0,1]  = "5verbs - retrained"
$row42[0,2]  = 56
$row42[0,3]  = 10
$row42[0,4]  = "e-4"
$row42[0,5]  = 100
$row42[0,6]  = 2
$row42[0,7]  = 8
$row42[0,8]  = 8
$row42[0,9]  = 0
$row42[0,10] = 0.22869999999999999
$row42[0,11] = 0.98899999999999999
$row42[0,12] = 0.48249999999999998
$row42[0,13] = 0.86960000000000004
$row42[0,14] = 0.29104477611940299
$ws.Range("A42:O42").Value = $row42
# D42 also carries the underline style seen elsewhere in this edit
$ws.Range("D42").Font.Underline = $true

$row43 = New-Object 'object[,]' 1,15
$row43[0,0]  = "IA-PUCP"
$row43[0,1]  = "5verbs - retrained"
$row43[0,2]  = 57
$row43[0,3]  = 10
$row43[0,4]  = "e-5"
$row43[0,5]  = 100
$row43[0,6]  = 2
$row43[0,7]  = 8
$row43[0,8]  = 8
$row43[0,9]  = 0
$row43[0,10] = 1.1679999999999999
$row43[0,11] = 0.70330000000000004
$row43[0,12] = 1.302
$row43[0,13] = 0.4783
$row43[0,14] = 0.56716417910447703
$ws.Range("A43:O43").Value = $row43

# Highlight the "Computer" column for the new retrained runs in yellow,
# same as was done for the earlier block of rows with the green fill.
$ws.Range("A40:A43").Interior.Color = 65535

# --- A couple of stray formatted-but-empty cells further down, carried
# over from the edit (rows 45-46) ---
$ws.Range("N45").Font.Underline = $true
$ws.Range("D46").Font.Underline = $true

# --- Column B needs to widen to fit the new, longer "ListWords" label ---
$ws.Columns("B").ColumnWidth = 16

# --- Leave the selection where the author left it when saving ---
$ws.Range("E45").Select()
